$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.08 = 7989.58 pesos`n✅ 7989.58 pesos = 2.08 = 937.14 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10/N12/O12 numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 480
$ws2.Range("O10").Value = 3835
$ws2.Range("N12").Value = 3845
$ws2.Range("O12").Value = 451
